# Update "想去人数" (number of people interested) counts for a few events.
# These updates apply identically to the "展览" sheet and the "全部类型"
# sheet, which mirrors the same data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 496
    $ws.Range("F6").Value = 682
    $ws.Range("F7").Value = 413
}
